$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all cells we will touch (to preserve the data as text, matching source)
# then set their literal text values.

$cellsToFormat = @(
    "D2",
    "E2",
    "G2",
    "D3",
    "E3",
    "G3",
    "D4",
    "E4",
    "G4",
    "D5",
    "E5",
    "G5",
    "D6",
    "E6",
    "G6",
    "D7",
    "E7",
    "G7",
    "D8",
    "E8",
    "G8",
    "E9",
    "G9",
    "D10",
    "E10",
    "G10",
    "D11",
    "E11",
    "G11",
    "D12",
    "E12",
    "G12",
    "D13",
    "E13",
    "G13",
    "D14",
    "E14",
    "G14",
    "D15",
    "E15",
    "G15",
    "D16",
    "E16",
    "G16",
    "D17",
    "E17",
    "G17",
    "E18",
    "G18",
    "E19",
    "G19",
    "D20",
    "E20",
    "G20",
    "D21",
    "E21",
    "G21",
    "E22",
    "G22",
    "D23",
    "E23",
    "G23",
    "D24",
    "E24",
    "G24",
    "D25",
    "E25",
    "G25",
    "D26",
    "E26",
    "G26",
    "G27",
    "G28",
    "G29",
    "G30",
    "G31",
    "G32",
    "G33",
    "G34",
    "G35",
    "G36",
    "G37",
    "G38",
    "D39",
    "E39",
    "G39",
    "D40",
    "E40",
    "G40",
    "D41",
    "E41",
    "G41",
    "D42",
    "E42",
    "G42",
    "D43",
    "E43",
    "G43",
    "D44",
    "E44",
    "G44",
    "D45",
    "E45",
    "G45",
    "D46",
    "E46",
    "G46",
    "E47",
    "G47",
    "D48",
    "G48",
    "E49",
    "G49",
    "E50",
    "G50",
    "D51",
    "E51",
    "G51"
)
foreach ($ref in $cellsToFormat) {
    $ws.Range($ref).NumberFormat = "@"
}

# Now assign the new literal values exactly as text
$ws.Range("D2").Value = "306.71"
$ws.Range("E2").Value = "-0.70%"
$ws.Range("G2").Value = "20"

$ws.Range("D3").Value = "38.81"
$ws.Range("E3").Value = "6.66%"
$ws.Range("G3").Value = "20"

$ws.Range("D4").Value = "5.109"
$ws.Range("E4").Value = "1.32%"
$ws.Range("G4").Value = "20"

$ws.Range("D5").Value = "0.08084"
$ws.Range("E5").Value = "-0.75%"
$ws.Range("G5").Value = "20"

$ws.Range("D6").Value = "1.948"
$ws.Range("E6").Value = "-7.13%"
$ws.Range("G6").Value = "20"

$ws.Range("D7").Value = "7.973"
$ws.Range("E7").Value = "1.48%"
$ws.Range("G7").Value = "20"

$ws.Range("D8").Value = "0.9309"
$ws.Range("E8").Value = "0.14%"
$ws.Range("G8").Value = "20"

$ws.Range("E9").Value = "2.09%"
$ws.Range("G9").Value = "20"

$ws.Range("D10").Value = "0.1933"
$ws.Range("E10").Value = "0.15%"
$ws.Range("G10").Value = "20"

$ws.Range("D11").Value = "0.09189"
$ws.Range("E11").Value = "0.86%"
$ws.Range("G11").Value = "20"

$ws.Range("D12").Value = "0.03509"
$ws.Range("E12").Value = "1.93%"
$ws.Range("G12").Value = "20"

$ws.Range("D13").Value = "0.09800"
$ws.Range("E13").Value = "-1.53%"
$ws.Range("G13").Value = "20"

$ws.Range("D14").Value = "0.001396"
$ws.Range("E14").Value = "-1.58%"
$ws.Range("G14").Value = "20"

$ws.Range("D15").Value = "0.005957"
$ws.Range("E15").Value = "-3.77%"
$ws.Range("G15").Value = "20"

$ws.Range("D16").Value = "3.782"
$ws.Range("E16").Value = "-1.38%"
$ws.Range("G16").Value = "20"

$ws.Range("D17").Value = "4.183"
$ws.Range("E17").Value = "0.73%"
$ws.Range("G17").Value = "20"

$ws.Range("E18").Value = "-0.86%"
$ws.Range("G18").Value = "20"

$ws.Range("E19").Value = "-0.16%"
$ws.Range("G19").Value = "20"

$ws.Range("D20").Value = "0.1303"
$ws.Range("E20").Value = "-1.11%"
$ws.Range("G20").Value = "20"

$ws.Range("D21").Value = "4.674"
$ws.Range("E21").Value = "-2.83%"
$ws.Range("G21").Value = "20"

$ws.Range("E22").Value = "3.31%"
$ws.Range("G22").Value = "20"

$ws.Range("D23").Value = "0.04380"
$ws.Range("E23").Value = "-0.07%"
$ws.Range("G23").Value = "20"

$ws.Range("D24").Value = "0.001237"
$ws.Range("E24").Value = "0.44%"
$ws.Range("G24").Value = "20"

$ws.Range("D25").Value = "0.004281"
$ws.Range("E25").Value = "1.87%"
$ws.Range("G25").Value = "20"

$ws.Range("D26").Value = "0.0001302"
$ws.Range("E26").Value = "0.36%"
$ws.Range("G26").Value = "20"

$ws.Range("G27").Value = "20"

$ws.Range("G28").Value = "20"

$ws.Range("G29").Value = "20"

$ws.Range("G30").Value = "20"

$ws.Range("G31").Value = "20"

$ws.Range("G32").Value = "20"

$ws.Range("G33").Value = "20"

$ws.Range("G34").Value = "20"

$ws.Range("G35").Value = "20"

$ws.Range("G36").Value = "20"

$ws.Range("G37").Value = "20"

$ws.Range("G38").Value = "20"

$ws.Range("D39").Value = "0.02040"
$ws.Range("E39").Value = "-0.09%"
$ws.Range("G39").Value = "20"

$ws.Range("D40").Value = "0.05080"
$ws.Range("E40").Value = "-1.57%"
$ws.Range("G40").Value = "20"

$ws.Range("D41").Value = "0.007523"
$ws.Range("E41").Value = "-0.21%"
$ws.Range("G41").Value = "20"

$ws.Range("D42").Value = "0.01027"
$ws.Range("E42").Value = "1.92%"
$ws.Range("G42").Value = "20"

$ws.Range("D43").Value = "0.1349"
$ws.Range("E43").Value = "-2.39%"
$ws.Range("G43").Value = "20"

$ws.Range("D44").Value = "0.002123"
$ws.Range("E44").Value = "-0.12%"
$ws.Range("G44").Value = "20"

$ws.Range("D45").Value = "0.009898"
$ws.Range("E45").Value = "1.61%"
$ws.Range("G45").Value = "20"

$ws.Range("D46").Value = "0.00006205"
$ws.Range("E46").Value = "-1.46%"
$ws.Range("G46").Value = "20"

$ws.Range("E47").Value = "0.30%"
$ws.Range("G47").Value = "20"

$ws.Range("D48").Value = "0.003109"
$ws.Range("G48").Value = "20"

$ws.Range("E49").Value = "0.28%"
$ws.Range("G49").Value = "20"

$ws.Range("E50").Value = "0.30%"
$ws.Range("G50").Value = "20"

$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "0.30%"
$ws.Range("G51").Value = "20"
